$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the updated symbol-list values (prices, coin links, labels) for the
# Sat Dec 24 2022 GitHub Actions data refresh.
# Column D holds numeric-looking price strings stored as text; force the
# cell's number format to Text ("@") before assigning so Excel does not
# coerce them into actual numbers (which would drop things like trailing
# zeros / leading zeros). Columns B, C and E hold non-numeric text, so a
# plain .Value assignment is sufficient for those.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '244.85'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.404'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06039'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8140'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9243'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1435'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07491'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.03394'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03047'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09428'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.011'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.04812'
$ws.Range("B17").Value = 'One'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0005941'
$ws.Range("E17").Value = '16OneONE'
$ws.Range("B18").Value = 'TigerCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.005681'
$ws.Range("E18").Value = '17TigerCashTCH'
$ws.Range("B19").Value = 'HotbitToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.004166'
$ws.Range("E19").Value = '18HotbitTokenHTB'
$ws.Range("B20").Value = 'BitKan'
$ws.Range("C20").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0009893'
$ws.Range("E20").Value = '19BitKanKAN'
$ws.Range("B21").Value = 'LEO'
$ws.Range("C21").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '3.665'
$ws.Range("E21").Value = '20LEOLEO'
$ws.Range("B22").Value = 'KuCoinToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.430'
$ws.Range("E22").Value = '21KuCoinTokenKCS'
$ws.Range("B23").Value = 'BTSEToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.183'
$ws.Range("E23").Value = '22BTSETokenBTSE'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0002901'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03999'
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1077'
$ws.Range("E41").Value = '40BKEXTokenBKK'
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.002720'
$ws.Range("E42").Value = '41CEJICEJI'
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.003052'
$ws.Range("E43").Value = '42KickTokenKICK'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.005791'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005234'
$ws.Range("E47").Value = '46CoinbaseStockTokenCOINBestin24h'
